# Update TrialsSetup 2026-02-20 12:00
#
# - Adds a new trial row (INNOVATE) to the Query1 table / Sheet1 data.
# - Corrects REMASTER (CLOU)'s "Days remaining" from 20 to 19.
# - Grows the table, its AutoFilter, the worksheet dimension and the
#   hidden ExternalData_1 defined name so they all cover the new row
#   (A1:C9 -> A1:C10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Days remaining" value for REMASTER (CLOU).
$ws.Range("B8").Value = 19

# Append the new trial as row 10 (Days remaining / Progress left blank,
# same as the other not-yet-started trials).
$ws.Range("A10").Value = "INNOVATE"

# Re-apply (explicit) General number formatting down column A for the
# data rows, including the freshly added row.
$ws.Range("A2:A10").NumberFormat = "General"

# Grow the query table (and its AutoFilter) to include the new row.
$table = $ws.ListObjects.Item("Query1")
$table.Resize($ws.Range("A1:C10"))

# Keep the hidden ExternalData_1 defined name (used by the Power Query
# connection) in sync with the new extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!ExternalData_1") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$10"
    }
}
